$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74; this shifts the existing rows 74-185 down to 75-186
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new translation-tracking entry
$ws.Range("A74").Value = "Initial & Final Surveillance Diagnosis"
$ws.Range("B74").Value = "TBT"
$ws.Range("C74").Value = "new"

# Fix the "Initial & Final Surveillance Diagnosis" row (now at 75): translate "&" to "and"
$ws.Range("A75").Value = "Initial and Final Surveillance Diagnosis"
$ws.Range("B75").Value = "រោគវិនិច្ឆ័យដំបូង និងចុងក្រោយនៃការចូលរួមការអង្កេតតាមដាន"

# Fix the "Susceptible & Intermediate ..." row (now at 145): translate "&" to "and"
$ws.Range("A145").Value = "Susceptible and Intermediate are always combined in this visualisation of co-resistances."
$ws.Range("B145").Value = "Susceptible and Intermediate​ គឺតែងតែរួមបញ្ជូលគ្នានៅក្នុង​គំនូសតាងនៃសហ-ភាពសុំា"
